# movement.xlsx edit: "만분율 float 변경, screen min,max 변경"
#  - B3 (the "moveSpeed" column's type row) changes from "int" to "float"
#  - the active/selected cell on Sheet1 moves from E9 to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 만분율(moveSpeed) 타입을 int -> float 으로 변경
$ws.Range("B3").Value = "float"

# 화면(선택 영역) 위치 변경
$ws.Range("B7").Select()
